# Add "Unallowed Application" check as a new row 27 in the Checklist sheet.
# Existing rows 27 ("Log browser URL") and 28 ("Workflow naming convention")
# shift down to 28 and 29 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 27 (shifts 27->28, 28->29).
$ws.Rows.Item(27).Insert()

# Copy the formatting (styles/borders/wrap) of row 26 ("Unallowed activity"),
# which uses the same column-style pattern (s=9 for A-D/F-G, s=6 for E) that
# the new row needs.
$ws.Range("A26:G26").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the explicit row height used by the author for this wrapped-text row.
$ws.Rows.Item(27).RowHeight = 130.5

# Fill in the values for the new check. Columns: A=Run, B=Issue,
# C=Check Filename, D=Arguments, E=Action, F=Explanation, G=Suggestion.
$ws.Range("A27").Value = "No"
$ws.Range("B27").Value = "Unallowed application"
$ws.Range("C27").Value = "Checks\Custom\UnallowedApplication.xaml"
$ws.Range("D27").Value = "{`n""WhiteList"" : """",`n""BlackList"": ""notepad.exe,explorer.exe""`n} "
$ws.Range("G27").Value = "Remove interactions with unallowed applications from the workflow or request the addition of the application to the whitelist (or its removal from the blacklist)."
$ws.Range("F27").Value = "Workflows should interact only with applications allowed by the Centre of Excellence (CoE). If an application is present in the blacklist or is not present in the whitelist of activities, it should not be used by the robot. The whitelist and the blacklist are passed as arguments to this check and contain names of applications' executable files. Each name must be specified as the application's executable file name (for example, ""notepad.exe"" instead of ""Notepad"")."
$ws.Range("E27").Value = "Fix"
